$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425; this shifts the existing row 425 (and
# everything below it) down by one, just like Excel's native "Insert Row".
$ws.Range("A425:T425").EntireRow.Insert()

# Populate the newly-inserted row 425 with the new record.
$ws.Range("A425").Value = 10
$ws.Range("B425").Value = "Vega Modelo de Temuco"
$ws.Range("C425").Value = "La Araucanía"
$ws.Range("D425").Value = 44722
$ws.Range("E425").Value = 9
$ws.Range("F425").Value = "Fruta"
$ws.Range("G425").Value = 100102
$ws.Range("H425").Value = "Cítricos"
$ws.Range("I425").Value = 100102004
$ws.Range("J425").Value = "Mandarina"
$ws.Range("K425").Value = "Clementina"
$ws.Range("L425").Value = "Primera"
$ws.Range("M425").Value = 285
$ws.Range("N425").Value = 15000
$ws.Range("O425").Value = 16000
$ws.Range("P425").Value = 15544
$ws.Range("Q425").Value = "$/bandeja 18 kilos"
$ws.Range("R425").Value = "Provincia de Limarí"
$ws.Range("S425").Value = 864
$ws.Range("T425").Value = 18
